# Applies the Thu Oct 26 11:21:56 UTC 2023 GitHub Actions cryptos-list refresh:
# updates Price/Volume(1h) figures row by row and fixes the swapped
# Cosmos/EthereumClassic (rows 26-27) and Kaspa/WEMIXToken (rows 46-47) pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.213.14'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.06%  '

# Row 3: Ethereum
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.828.55'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.48%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.17'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.23%  '

# Row 6: XRP
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.559'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.26%  '

# Row 7: USDC
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.06%  '

# Row 8: Solana
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.00'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.44%  '

# Row 9: Cardano
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.84%  '

# Row 10: Dogecoin
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0727'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +10.76%  '

# Row 11: TRON
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0930'

# Row 12: WrappedliquidstakedEther2.0
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.092.04'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.59%  '

# Row 13: WrappedEther
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.824.69'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.35%  '

# Row 14: Chainlink
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.86'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.66%  '

# Row 15: Polygon
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.645'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.92%  '

# Row 16: WrappedBTC
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.259.73'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.15%  '

# Row 17: Polkadot
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.37%  '

# Row 18: Litecoin
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.82'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.55%  '

# Row 19: BitcoinCash
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '251.63'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.14%  '

# Row 20: ShibaInu
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0795'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +7.58%  '

# Row 21: Avalanche
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.20'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +8.15%  '

# Row 22: Dai
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.04%  '

# Row 23: Uniswap
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.01%  '

# Row 24: Toncoin
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.40%  '

# Row 25: Monero
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '160.68'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.11%  '

# Row 26: Cosmos
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.72'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.95%  '

# Row 27: EthereumClassic
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.29'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.91%  '

# Row 28: Stellar
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.42%  '

# Row 29: BinanceUSD
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.06%  '

# Row 30: Hedera
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.15%  '

# Row 31: Filecoin
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.79'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.54%  '

# Row 32: PancakeSwap
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.27%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.59'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.17%  '

# Row 34: LidoDAOToken
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.84%  '

# Row 35: Maker
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.441.45'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.10%  '

# Row 36: ImmutableX
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.75%  '

# Row 37: TrustWalletToken
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.43%  '

# Row 38: VeChain
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.48%  '

# Row 39: ARBITRUM
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.967'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +8.56%  '

# Row 40: Aave
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '81.93'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.29%  '

# Row 41: MXToken
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.76'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.28%  '

# Row 42: HuobiToken
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.34'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.07%  '

# Row 43: RenderToken
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +5.45%  '

# Row 44: FraxShare
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.11'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.77%  '

# Row 45: RocketPoolETH
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.989.33'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.46%  '

# Row 46: Kaspa
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.06'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.61%  '

# Row 47: WEMIXToken
$ws.Range('B47').Value = 'Kaspa'
$ws.Range('C47').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0498'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.54%  '

# Row 48: Quant
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '107.67'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +9.30%  '

# Row 49: PaxDollar
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.02%  '

# Row 50: InjectiveProtocol
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.12%  '

# Row 51: BabyDogeCoin
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0124'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.56%  '
